# Daily crypto price/volume refresh (GitHub Actions data pull).
#
# Column D ("Price") holds plain text even when the text looks numeric
# (e.g. "680.77"), so force text formatting before writing any value that
# would otherwise be auto-converted to a float by Excel's smart-typing.
# Column E ("Volume(1h)") values already contain surrounding whitespace and
# a trailing "%", so they are written as text without any extra handling.
#
# Rows 33/34 additionally swap coin identity: Binance-Peg BSC-USD now
# outranks EthereumClassic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.340.12"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.686.66"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "679.29"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.03"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.442"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").Value = "4.309.47"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "3.676.58"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "69.298.76"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.36"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.88"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.99"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "3.834.69"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -5.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.89"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.15"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.69"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.91"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "3.677.01"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.26"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0904"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "170.57"
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.43"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.21"
$ws.Range("E46").Value = "  -5.41%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  -2.99%  "
